# "updated main GSC export data"
#
# Appends the next day (2025-12-11) as a new row at the bottom of the
# "Chart" sheet's A:C data table (Date / Non-HTTPS URLs / HTTPS URLs),
# mirroring every previous daily row already in the sheet.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

$newRow = 67  # one past the last existing data row (66)

# Column A holds dates as plain text (e.g. "2025-12-10"), not real date
# serials. Assigning a date-shaped string straight to .Value makes Excel
# auto-recognise it as a date and reformat the cell, which we don't want
# here. Routing it through a formula result + "paste values" keeps it a
# plain text value, just like the rest of column A.
$stage = $chart.Range("Z1")
$stage.Formula = "=""2025-12-11"""
$stage.Copy()
$chart.Cells.Item($newRow, 1).PasteSpecial(-4163)  # xlPasteValues
$stage.Clear()

$chart.Cells.Item($newRow, 2).Value = 0
$chart.Cells.Item($newRow, 3).Value = 29
